# Replace the synthetic "long" dataset on Sheet1 with the regenerated
# multi-index / multi-table sample (weekly Date x Channel x Metric x Value
# rows) produced by the updated dataset-generation script.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46050, "Radio", "GRPs", 9),
    @(46078, "TV", "Spend", 57),
    @(46050, "Radio", "Spend", 139),
    @(46078, "Radio", "GRPs", 5),
    @(46050, "TV", "Spend", 99),
    @(46071, "TV", "GRPs", 9),
    @(46029, "TV", "Spend", 114),
    @(46036, "Radio", "GRPs", 4),
    @(46015, "TV", "Spend", 134),
    @(46050, "TV", "GRPs", 1),
    @(46022, "TV", "Spend", 113),
    @(46008, "Radio", "GRPs", 3),
    @(46015, "Radio", "Spend", 69),
    @(46071, "Radio", "GRPs", 10),
    @(46064, "TV", "GRPs", 3),
    @(46029, "TV", "GRPs", 10),
    @(46036, "TV", "GRPs", 9),
    @(46085, "TV", "GRPs", 3),
    @(46029, "Radio", "GRPs", 2),
    @(46022, "TV", "GRPs", 8),
    @(46036, "TV", "Spend", 197),
    @(46057, "TV", "GRPs", 1),
    @(46064, "TV", "Spend", 105),
    @(46085, "Radio", "Spend", 99),
    @(46008, "TV", "Spend", 96),
    @(46057, "Radio", "GRPs", 8),
    @(46057, "Radio", "Spend", 126),
    @(46015, "TV", "GRPs", 1),
    @(46022, "Radio", "GRPs", 1),
    @(46064, "Radio", "Spend", 94),
    @(46008, "TV", "GRPs", 3),
    @(46085, "TV", "Spend", 53),
    @(46043, "TV", "GRPs", 7),
    @(46008, "Radio", "Spend", 89),
    @(46071, "TV", "Spend", 51),
    @(46064, "Radio", "GRPs", 9),
    @(46043, "Radio", "GRPs", 8),
    @(46078, "TV", "GRPs", 4),
    @(46022, "Radio", "Spend", 115),
    @(46043, "Radio", "Spend", 73),
    @(46036, "Radio", "Spend", 151),
    @(46085, "Radio", "GRPs", 10),
    @(46043, "TV", "Spend", 98),
    @(46015, "Radio", "GRPs", 1),
    @(46078, "Radio", "Spend", 126),
    @(46029, "Radio", "Spend", 157),
    @(46057, "TV", "Spend", 191),
    @(46071, "Radio", "Spend", 171)
)

$numRows = $data.Count
$numCols = 4

# Build a 2-D array so the whole block can be written in a single
# Range.Value2 assignment (row order above = Date serial, Channel, Metric, Value).
$arr = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$rng.Value2 = $arr

# Column A holds date serials - keep the existing date/time display format
# (same numFmt already used by the header rows' style) for every data row,
# including the newly added ones (rows 33-49).
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
